$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.805.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.11%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.350.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.30%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'544.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.79%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -3.23%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.349.47"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +1.84%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'5.32"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.49%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.344"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'24.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.88%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'2.773.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.19%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'60.710.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.01%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -1.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.348.35"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.36%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'10.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.28%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'Polkadot"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.51%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "'BitcoinCash"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").Value = "'319.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.89%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.18%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.03%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'63.35"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.58%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.69"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -6.62%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'8.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +8.19%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B28").Value = "'InternetComputer(DFINITY)"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'7.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.25%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'Bittensor"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'498.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -3.76%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Fetch.AI"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.87%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value = "'PEPE"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'0.0₃0862"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -6.76%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "'Kaspa"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'0.145"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +0.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "'PancakeSwap"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.79"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.04%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'ImmutableX"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'1.50"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.44%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'FirstDigitalUSD"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'NEARProtocol"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'4.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.97%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'PolygonEcosystemToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.376"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.47%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'EthereumClassic"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'18.47"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.20%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Stacks"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +5.94%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'RenderToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'5.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.76%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'Monero"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'144.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.28%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'USDe"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.05%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'OKB"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'40.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.91%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = "'Aave"
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = "'143.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.09%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("B45").Value = "'Filecoin"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.56"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.50%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'dogwifhat"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'2.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -8.44%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Hedera"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0517"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.25%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "'InjectiveProtocol"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'19.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.87%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'Mantle"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'0.570"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.12%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'Stellar"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0901"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.13%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'VeChain"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.0221"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.77%  "
$ws.Range("E51").Style = "Normal"
